$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha Ponto")

# Row 8 (new timesheet entry: 13:00 - 13:30, SITS, "Implemented pico font")
# B8/C8/D8 don't yet carry any formatting, so copy the date/time number
# formats from the row above before writing the values.
$ws.Cells.Item(2, 2).Copy()
$ws.Cells.Item(8, 2).PasteSpecial(-4122)
$ws.Cells.Item(2, 3).Copy()
$ws.Cells.Item(8, 3).PasteSpecial(-4122)
$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4122)

$ws.Cells.Item(8, 2).Value2 = 44902
$ws.Cells.Item(8, 3).Value2 = 0.54166666666666663
$ws.Cells.Item(8, 4).Value2 = 0.5625
$ws.Cells.Item(8, 6).Value2 = "SITS"
$ws.Cells.Item(8, 7).Value2 = "Implemented pico font"

# Row 9 (new timesheet entry: 13:30 - 13:45, SITS, "Implemented coin counter display")
# B9/C9/D9 already carry the right number formats, so just set the values.
$ws.Cells.Item(9, 2).Value2 = 44902
$ws.Cells.Item(9, 3).Value2 = 0.5625
$ws.Cells.Item(9, 4).Value2 = 0.57291666666666663
$ws.Cells.Item(9, 6).Value2 = "SITS"
$ws.Cells.Item(9, 7).Value2 = "Implemented coin counter display"

# Move the active selection as recorded in the saved view state.
[void]$ws.Range("G12").Select()
